$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function OrdinalCompare($a, $b) {
    $la = $a.Length
    $lb = $b.Length
    $cn = [Math]::Min($la, $lb)
    for ($ci = 0; $ci -lt $cn; $ci++) {
        $ca = [int][char]$a[$ci]
        $cb = [int][char]$b[$ci]
        if ($ca -lt $cb) { return -1 }
        if ($ca -gt $cb) { return 1 }
    }
    if ($la -lt $lb) { return -1 }
    if ($la -gt $lb) { return 1 }
    return 0
}

function OrdinalSort($arr) {
    $sn = $arr.Length
    for ($si = 1; $si -lt $sn; $si++) {
        $key = $arr[$si]
        $sj = $si - 1
        $cmp = OrdinalCompare $arr[$sj] $key
        while ($sj -ge 0 -and $cmp -gt 0) {
            $arr[$sj + 1] = $arr[$sj]
            $sj = $sj - 1
            if ($sj -ge 0) {
                $cmp = OrdinalCompare $arr[$sj] $key
            }
        }
        $arr[$sj + 1] = $key
    }
    return $arr
}

$lastRow = $ws.UsedRange.Rows.Count
Write-Host "lastRow=$lastRow"

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $orig = $cell.Text
    if ($orig -ne "") {
        $parts = $orig.Split(",")
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }
        if ($trimmed.Length -gt 1) {
            $sortedParts = OrdinalSort $trimmed
            $newVal = $sortedParts -join ", "
            if ($newVal -ne $orig) {
                $cell.Value = $newVal
            }
        }
    }
}

Write-Host "Done"
